$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 16: update the table's style id on the graphicFrame shape.
# ---------------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$tblShape = $null
for ($i = 1; $i -le $s16.Shapes.Count; $i++) {
    $cand = $s16.Shapes.Item($i)
    if ($cand.HasTable) {
        $tblShape = $cand
        break
    }
}
$tbl = $tblShape.Table
$tbl.ApplyStyle("{9F32882F-94D8-47E9-A19B-EF3A0CF9EB22}")

# ---------------------------------------------------------------------------
# 2) Swap the theme colour scheme: the deck's active theme (the one used by
#    the slide master) moves from the "Integral" palette to the classic
#    "Office Theme" palette.
# ---------------------------------------------------------------------------
function HexToRgbInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToRgbInt($officeColors[$i - 1])
}
